# Journal de bord - notes du TP1 pour TPA1 et TPA3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Fill in the data for rows 4 and 5 (new TP1 session entries)
#    Order of entry matters: it drives the shared-string table order.
# ------------------------------------------------------------------

# Row 4: date, teacher, seance type, group A3 = X, description
$ws.Range("B4").Value = "MPAL"
$ws.Range("G4").Value = "Lecture active du CDC, mise en commun en binome. 2 questions Elaastic."

# Row 5: date, teacher, seance type, group A1 = X, description
$ws.Range("G5").Value = "Lecture active du CDC, mise en commun en binome. 3 questions Elaastic."

# Re-use existing shared strings ("TP", "X") for the remaining known cells
$ws.Range("C4").Value = "TP"
$ws.Range("F4").Value = "X"

$ws.Range("B5").Value = "MPAL"
$ws.Range("C5").Value = "TP"
$ws.Range("D5").Value = "X"

# Dates (12 sept. 2023) - copy format from A2 (already date-formatted) then set value
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 45181

$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 45181

# ------------------------------------------------------------------
# 2) Add the new "Commentaires" table column (9th column)
# ------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Add()
$col.Range.Cells(1).Value = "Commentaires"

# Body cells for the new column: give every data row the same base
# style as the rest of the table (copy formats from an existing cell)
$ws.Range("A3").Copy()
$ws.Range("I2:I3").PasteSpecial(-4122)
$ws.Range("I6:I31").PasteSpecial(-4122)

# Rows 4 & 5 comments use the wrap-text style already used by column H
$ws.Range("H2").Copy()
$ws.Range("I4:I5").PasteSpecial(-4122)

$ws.Range("I4").Value = "Ce qui est ressorti des discussion : pas d'exigences liées à l'ergonomie ou l'architecture hierarchique du site. Clarification de l'obtention du niveau de privilège. Mise en commun avec la classe et les très bonnes questions soulevées ont amené à devancer un peu les questions Elaastic...`nElaastic 1 : Peu de réponses liées au caractère évolutif et imprécis du CDC, qu'il faudrait transformer en US. Plutôt tendance à dire qu'il faut modifier le CDC.`nElaastic 2 : Peu efficace. Surtout des réponses liées au fait qu'il ne s'agissait pas d'une fonctionnalité mais d'une notion de compatibilité. Personne n'a soulevé la notion de fonctionnalité/exigence transversale et non finissable."

$ws.Range("I5").Value = "Discussions : beaucoup de questionnements. `nCertains ne voient pas ce qui cloche : ""c'est normal de se poser des questions, on y répondra au fur et à mesure avec le client car méthode Agile"")`nPas de mise en commun des questionnements avant Elaastic.`nElaastic 1 : tout le monde est d'accord pour dire que le CDC est imprécis. Les meilleurs notes disent ce qu'il manque. Peu disent que c'est normal et que la méthode Agile servira à clarifier. `nElaastic 2 : bon repérage du caractère non fonctionnel, mais peu disent la notion de non finissable. Certain confondent exigence qualité, technique, accessibilité...`nElaastic 3 : remarque générale sur la durée de développement, mais pas sur la couverture fonctionnelle de 100%."

# ------------------------------------------------------------------
# 3) Row heights for the affected rows
# ------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 46.5
$ws.Rows.Item(3).RowHeight = 170.5
$ws.Rows.Item(4).RowHeight = 139.5
$ws.Rows.Item(5).RowHeight = 155

# ------------------------------------------------------------------
# 4) Column width for the new "Commentaires" column
# ------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 78.92

# ------------------------------------------------------------------
# 5) Selection / active cell
# ------------------------------------------------------------------
$ws.Range("I6").Select()
